# Sync automático del tracker - marcar predicciones como Completadas
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamp = "2025-08-22 21:17:28"

$updates = @(
    @{ Row = 30; Result = "Home Win"; Real = "Fallo";   Profit = -6.2; ROI = -100 },
    @{ Row = 32; Result = "Home Win"; Real = "Fallo";   Profit = -2.8; ROI = -100 },
    @{ Row = 34; Result = "Away Win"; Real = "Fallo";   Profit = -6;   ROI = -100 },
    @{ Row = 35; Result = "Draw";     Real = "Fallo";   Profit = -5.3; ROI = -100 },
    @{ Row = 36; Result = "Draw";     Real = "Fallo";   Profit = -2.7; ROI = -100 },
    @{ Row = 37; Result = "Home Win"; Real = "Fallo";   Profit = -1.7; ROI = -100 },
    @{ Row = 38; Result = "Home Win"; Real = "Fallo";   Profit = -1.8; ROI = -100 },
    @{ Row = 39; Result = "Away Win"; Real = "Fallo";   Profit = -4.8; ROI = -100 },
    @{ Row = 42; Result = "Home Win"; Real = "Fallo";   Profit = -2.4; ROI = -100 },
    @{ Row = 50; Result = "Draw";     Real = "Fallo";   Profit = -4.5; ROI = -100 },
    @{ Row = 86; Result = "Draw";     Real = "Fallo";   Profit = -4.5; ROI = -100 },
    @{ Row = 94; Result = "Away Win"; Real = "Acierto"; Profit = 8.99; ROI = 333 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Range("L$r").Value = "Completed"
    $ws.Range("M$r").Value = $u.Result
    $ws.Range("N$r").Value = $u.Real
    $ws.Range("O$r").Value = $u.Profit
    $ws.Range("P$r").Value = $u.ROI
    $ws.Range("Q$r").Value = $timestamp
}
